$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Header / summary block updates
# -----------------------------------------------------------------
$ws.Range("E11").Value = 1148000
$ws.Range("C13").Value = 6
$ws.Range("F13").Value = 4

# Widen column D slightly (OOXML raw width 33.08984375 -> 35)
$ws.Range("D1").ColumnWidth = 34.24015625

# -----------------------------------------------------------------
# Expand the worker/period detail table from 6 rows (16-21) to
# 20 rows (16-35). Insert 14 blank rows right after the existing
# block (after row 21), which pushes the old signature rows
# (26/27) down to (40/41).
# -----------------------------------------------------------------
$ws.Range("22:35").Insert() | Out-Null

# Re-apply the "normal" data-row formatting (as seen on row 16) to
# the freshly inserted rows 22-34, and the "final row" formatting
# (as seen on the old last row, 21) to the new final row, 35.
$ws.Range("B16:J16").Copy() | Out-Null
$ws.Range("B22:J34").PasteSpecial(-4122) | Out-Null
$ws.Range("B21:J21").Copy() | Out-Null
$ws.Range("B35:J35").PasteSpecial(-4122) | Out-Null

# -----------------------------------------------------------------
# Write the new detail-table contents (7 workers, 20 rows total).
# -----------------------------------------------------------------
$data = @(
    @(16, "CC", "1005581984", "STEFANNY SUAREZ CONTRERAS",     "2507", 64000, 1600000),
    @(17, "CC", "1005581984", "STEFANNY SUAREZ CONTRERAS",     "2506", 64000, 1600000),
    @(18, "CC", "1005581984", "STEFANNY SUAREZ CONTRERAS",     "2505", 64000, 1600000),
    @(19, "CC", "1091659312", "JUAN CARLOS GUEVARA CALDERON",  "2507", 52000, 1300000),
    @(20, "CC", "1091659312", "JUAN CARLOS GUEVARA CALDERON",  "2506", 52000, 1300000),
    @(21, "CC", "1091659312", "JUAN CARLOS GUEVARA CALDERON",  "2505", 52000, 1300000),
    @(22, "CC", "1091659312", "JUAN CARLOS GUEVARA CALDERON",  "2504", 52000, 1300000),
    @(23, "CC", "18925242",   "JORGE HUMBERTO CALDERON YEPES", "2507", 52000, 1300000),
    @(24, "CC", "18925242",   "JORGE HUMBERTO CALDERON YEPES", "2506", 52000, 1300000),
    @(25, "CC", "18925242",   "JORGE HUMBERTO CALDERON YEPES", "2505", 52000, 1300000),
    @(26, "CC", "1016092355", "JORGE IVAN CUADRO PALOMINO",    "2507", 76000, 1900000),
    @(27, "CC", "1016092355", "JORGE IVAN CUADRO PALOMINO",    "2506", 76000, 1900000),
    @(28, "CC", "1016092355", "JORGE IVAN CUADRO PALOMINO",    "2505", 76000, 1900000),
    @(29, "CC", "12459249",   "ISMAEL CONTRERAS BECERRA",      "2507", 52000, 1300000),
    @(30, "CC", "12459249",   "ISMAEL CONTRERAS BECERRA",      "2506", 52000, 1300000),
    @(31, "CC", "12459249",   "ISMAEL CONTRERAS BECERRA",      "2505", 52000, 1300000),
    @(32, "CC", "12459249",   "ISMAEL CONTRERAS BECERRA",      "2504", 52000, 1300000),
    @(33, "CC", "92191463",   "ENILSO MANUEL PINEDA BETIN",    "2507", 52000, 1300000),
    @(34, "CC", "92191463",   "ENILSO MANUEL PINEDA BETIN",    "2506", 52000, 1300000),
    @(35, "CC", "92191463",   "ENILSO MANUEL PINEDA BETIN",    "2505", 52000, 1300000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}
